$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$rows = @(
    @(2, 44.15999984741211, 43.9900016784668, 46, 42.52000045776367, 94500000),
    @(3, 37.2599983215332, 37.11000061035156, 37.52999877929688, 35.11000061035156, 94500000),
    @(4, 28.97999954223633, 28.11000061035156, 30.27000045776367, 27.56999969482422, 94500000),
    @(5, 25.97999954223633, 24.34000015258789, 26.35000038146973, 21.97999954223633, 94500000),
    @(6, 25.81999969482422, 25.29999923706055, 26.71999931335449, 24.36000061035156, 94500000),
    @(7, 25.21999931335449, 28.3799991607666, 28.75, 24.51000022888184, 94500000),
    @(8, 30.07999992370605, 26.95999908447266, 30.77000045776367, 26.63999938964844, 94500000),
    @(9, 27.32999992370605, 29.36000061035156, 29.45999908447266, 27.13999938964844, 94500000),
    @(10, 31.09000015258789, 29.18000030517578, 32.20000076293945, 28.68000030517578, 94500000),
    @(11, 29.55999946594238, 31.81999969482422, 32.79999923706055, 28.32999992370605, 94500000),
    @(12, 33.79999923706055, 33.45000076293945, 34.29999923706055, 31.77000045776367, 94500000),
    @(13, 38.79999923706055, 40.5, 41.31999969482422, 37.65000152587891, 94500000),
    @(14, 39.72000122070312, 40.91999816894531, 41.90000152587891, 38.91999816894531, 94500000),
    @(15, 40.02000045776367, 38.29000091552734, 43.40999984741211, 37.29000091552734, 94500000),
    @(16, 37.93999862670898, 36.40000152587891, 38.08000183105469, 34.06000137329102, 94500000),
    @(17, 37.65000152587891, 44.38000106811523, 44.58000183105469, 36.47999954223633, 94500000),
    @(18, 44.08000183105469, 45.47000122070312, 46.15999984741211, 42.86999893188477, 94500000),
    @(19, 36.40999984741211, 36.61999893188477, 37.61999893188477, 34.45000076293945, 94500000),
    @(20, 31.23999977111816, 29.93000030517578, 31.26000022888184, 27.95000076293945, 94500000),
    @(21, 26.98999977111816, 24.34000015258789, 27.79000091552734, 24.14999961853028, 94500000),
    @(22, 19.76000022888184, 24.59000015258789, 25.65999984741211, 19, 94500000),
    @(23, 20.86000061035156, 21, 21.3700008392334, 19.01000022888184, 94500000),
    @(24, 22.65999984741211, 18.3700008392334, 22.95999908447266, 17.98999977111816, 94500000),
    @(25, 22.46999931335449, 26.89999961853028, 33.79999923706055, 21.95999908447266, 94500000),
    @(26, 39.40000152587891, 49.47000122070312, 54.88000106811523, 37.84999847412109, 94500000),
    @(27, 49.81999969482422, 49.65999984741211, 50.95999908447266, 43.02000045776367, 94500000),
    @(28, 57.68999862670898, 56.56000137329102, 59.38000106811523, 54.11000061035156, 94500000),
    @(29, 42.81000137329102, 40.34000015258789, 46.2599983215332, 37.34999847412109, 94500000),
    @(30, 49.29999923706055, 41.34999847412109, 50, 41.22999954223633, 94500000),
    @(31, 37.06000137329102, 38.29000091552734, 39.84000015258789, 36.18000030517578, 94500000),
    @(32, 31.23999977111816, 31.59000015258789, 32.36999893188477, 28.64999961853028, 94500000),
    @(33, 34.16999816894531, 34.88000106811523, 35.04000091552734, 32.33000183105469, 94500000),
    @(34, 40.15000152587891, 38.70999908447266, 41.09999847412109, 37.66999816894531, 94500000),
    @(35, 53.29999923706055, 56.84999847412109, 57.20000076293945, 51.91999816894531, 94500000),
    @(36, 45.0099983215332, 42.72000122070312, 46.72999954223633, 41.84000015258789, 94500000),
    @(37, 43.20999908447266, 46.18000030517578, 49.29000091552734, 42.22000122070312, 94500000),
    @(38, 38.81000137329102, 37.09999847412109, 39.11999893188477, 35.36000061035156, 94500000),
    @(39, 34.54000091552734, 32.41999816894531, 35.59999847412109, 30.52000045776367, 94500000),
    @(40, 30.31999969482422, 32.22999954223633, 33.13999938964844, 29.42000007629395, 94500000),
    @(41, 31.52000045776367, 31.90999984741211, 33.02999877929688, 29.77000045776367, 94500000),
    @(42, 22.54999923706055, 21.5, 23.19000053405762, 18.43000030517578, 94500000),
    @(43, 22.5, 20.93000030517578, 23.57999992370605, 20.88999938964844, 94500000)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 4).Value = $r[1]
    $ws.Cells.Item($rowNum, 5).Value = $r[2]
    $ws.Cells.Item($rowNum, 6).Value = $r[3]
    $ws.Cells.Item($rowNum, 7).Value = $r[4]
    $ws.Cells.Item($rowNum, 8).Value = $r[5]
    $ws.Cells.Item($rowNum, 9).Value = "TDC"
}